$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("A2:F33")
$key = $ws.Range("B2:B33")

$sortRange.Sort($key, 2)

$highlightCompanies = @("• Guj Alkali", "• Century PlyBoard", "• CCL Products", "• Cera Sanitary")

for ($r = 2; $r -le 33; $r++) {
    $name = $ws.Cells.Item($r, 1).Value
    if ($highlightCompanies -contains $name) {
        $ws.Range("A" + $r + ":F" + $r).Interior.Color = 65535
    }
}

$ws.Range("A30").Select()
